$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Append " TBD." to the "Contexts parsing ..." paragraph.
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -clike "*Recognize context types from (surrounding) reified kind types / rules (link grammar).*") {
        $r = $p.Range
        # Trim the trailing paragraph mark so we append right after the
        # existing sentence, inside the same run.
        $r.End = $r.End - 1
        $r.InsertAfter(" TBD.")
        break
    }
}

# ------------------------------------------------------------------
# 2) After the "Lower hierarchy layers contexts mapping ... TBD."
#    paragraph, insert a new "Contexts (TBD):" block:
#
#      (blank)
#      Contexts (TBD):
#      (blank)
#      Data layer: Resource, Statement, Role, Kind, Relation contexts (aggregation).
#      (blank)
#      Information layer (occurrences / interactions): Data layer contexts statements products as new contexts (aggregation).
#      (blank)
#      Knowledge layer (Dimensional / DCI Contexts): Information contexts statements products as new contexts (aggregation).
# ------------------------------------------------------------------
$anchorIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -clike "*Lower hierarchy layers contexts mapping*Roles to Kinds*") {
        $anchorIndex = $i
        break
    }
}

$newTexts = @(
    "",
    "Contexts (TBD):",
    "",
    "Data layer: Resource, Statement, Role, Kind, Relation contexts (aggregation).",
    "",
    "Information layer (occurrences / interactions): Data layer contexts statements products as new contexts (aggregation).",
    "",
    "Knowledge layer (Dimensional / DCI Contexts): Information contexts statements products as new contexts (aggregation)."
)

$insertAt = $anchorIndex
foreach ($text in $newTexts) {
    $cur = $d.Paragraphs.Item($insertAt)
    $cur.Range.InsertParagraphAfter()
    $insertAt = $insertAt + 1
    if ($text -ne "") {
        $newP = $d.Paragraphs.Item($insertAt)
        $newP.Range.InsertBefore($text)
    }
}
